$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.435.33'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.725.51'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.07'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4797'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2676'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06222'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.730.04'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.68'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6153'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.542'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.13'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.0000'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.458.24'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006936'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.70'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.953.02'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.538'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.902'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.311'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.35'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.39'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.793'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.407'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.74'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.971'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08037'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.732'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04547'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.620'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6372'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9880'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9386'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.985'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +4.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.413'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '107.27'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.30%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01497'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.635'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +10.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3907'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.987'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +12.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1191'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05318'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.917'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.33%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.91'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.265'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3417'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.11%  '
